# Update automàtic: dades i banners [2026-02-10 20:20]
# Refresh the per-station weather snapshot (extraction timestamp + measured
# values) in the Dades_Meteo sheet to match the 20:20 automated run.
# Note: column H holds humidity values formatted as plain text like "90%".
# Prefixing the literal value with a leading apostrophe forces Excel to
# store it as text (quote-prefixed) instead of auto-converting it into a
# numeric percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

$ws.Range('E2').Value = '2026-02-10 20:18:25'
$ws.Range('I2').Value = '46.7 mm'
$ws.Range('O2').Value = '0.9 °C'
$ws.Range('E3').Value = '2026-02-10 20:18:27'
$ws.Range('I3').Value = '25.3 mm'
$ws.Range('O3').Value = '0.8 °C'
$ws.Range('E4').Value = '2026-02-10 20:18:30'
$ws.Range('J4').Value = '1003.8 hPa'
$ws.Range('E5').Value = '2026-02-10 20:18:32'
$ws.Range('I5').Value = '34.0 mm'
$ws.Range('E6').Value = '2026-02-10 20:18:34'
$ws.Range('L6').Value = '13.3 km/h - 347º 19:52 TU'
$ws.Range('E7').Value = '2026-02-10 20:18:37'
$ws.Range('J7').Value = '1004.5 hPa'
$ws.Range('O7').Value = '15.3 °C'
$ws.Range('E8').Value = '2026-02-10 20:18:39'
$ws.Range('J8').Value = '1004.4 hPa'
$ws.Range('O8').Value = '12.1 °C'
$ws.Range('E9').Value = '2026-02-10 20:18:42'
$ws.Range('I9').Value = '0.7 mm'
$ws.Range('O9').Value = '8.9 °C'
$ws.Range('E10').Value = '2026-02-10 20:18:44'
$ws.Range('H10').Value = '''90%'
$ws.Range('E11').Value = '2026-02-10 20:18:47'
$ws.Range('I11').Value = '2.6 mm'
$ws.Range('E12').Value = '2026-02-10 20:18:49'
$ws.Range('I12').Value = '0.6 mm'
$ws.Range('O12').Value = '9.1 °C'
$ws.Range('E13').Value = '2026-02-10 20:18:51'
$ws.Range('I13').Value = '9.3 mm'
$ws.Range('E14').Value = '2026-02-10 20:18:54'
$ws.Range('E15').Value = '2026-02-10 20:18:56'
$ws.Range('I15').Value = '0.9 mm'
$ws.Range('E16').Value = '2026-02-10 20:18:59'
$ws.Range('I16').Value = '25.5 mm'
$ws.Range('E17').Value = '2026-02-10 20:19:01'
$ws.Range('E18').Value = '2026-02-10 20:19:03'
$ws.Range('E19').Value = '2026-02-10 20:19:06'
$ws.Range('O19').Value = '6.7 °C'
$ws.Range('E20').Value = '2026-02-10 20:19:08'
$ws.Range('I20').Value = '8.4 mm'
$ws.Range('O20').Value = '0.5 °C'
$ws.Range('E21').Value = '2026-02-10 20:19:10'
$ws.Range('H21').Value = '''91%'
$ws.Range('J21').Value = '1006.1 hPa'
$ws.Range('O21').Value = '7.3 °C'
$ws.Range('E22').Value = '2026-02-10 20:19:13'
$ws.Range('I22').Value = '9.3 mm'
$ws.Range('E23').Value = '2026-02-10 20:19:15'
$ws.Range('I23').Value = '27.0 mm'
$ws.Range('E24').Value = '2026-02-10 20:19:18'
$ws.Range('H24').Value = '''94%'
$ws.Range('E25').Value = '2026-02-10 20:19:20'
$ws.Range('I25').Value = '19.6 mm'
$ws.Range('E26').Value = '2026-02-10 20:19:23'
$ws.Range('I26').Value = '0.5 mm'
$ws.Range('E27').Value = '2026-02-10 20:19:25'
$ws.Range('I27').Value = '11.9 mm'
$ws.Range('E28').Value = '2026-02-10 20:19:28'
$ws.Range('I28').Value = '1.3 mm'
$ws.Range('O28').Value = '9.0 °C'
$ws.Range('E29').Value = '2026-02-10 20:19:30'
$ws.Range('O29').Value = '10.7 °C'
$ws.Range('E30').Value = '2026-02-10 20:19:33'
$ws.Range('I30').Value = '0.5 mm'
$ws.Range('O30').Value = '9.3 °C'
$ws.Range('E31').Value = '2026-02-10 20:19:35'
$ws.Range('I31').Value = '1.2 mm'
$ws.Range('M31').Value = '15.2 °C 19:36 TU'
$ws.Range('O31').Value = '10.4 °C'
$ws.Range('E32').Value = '2026-02-10 20:19:38'
$ws.Range('E33').Value = '2026-02-10 20:19:40'
$ws.Range('I33').Value = '11.0 mm'
$ws.Range('J33').Value = '1006.4 hPa'
$ws.Range('L33').Value = '11.9 km/h - 325º 19:44 TU'
$ws.Range('O33').Value = '4.3 °C'
$ws.Range('E34').Value = '2026-02-10 20:19:43'
$ws.Range('H34').Value = '''79%'
$ws.Range('I34').Value = '14.2 mm'
$ws.Range('E35').Value = '2026-02-10 20:19:45'
$ws.Range('E36').Value = '2026-02-10 20:19:47'
$ws.Range('I36').Value = '0.9 mm'
$ws.Range('E37').Value = '2026-02-10 20:19:50'
$ws.Range('H37').Value = '''91%'
$ws.Range('I37').Value = '3.3 mm'
$ws.Range('J37').Value = '1005.6 hPa'
$ws.Range('O37').Value = '6.7 °C'
$ws.Range('E38').Value = '2026-02-10 20:19:52'
$ws.Range('O38').Value = '11.0 °C'
$ws.Range('E39').Value = '2026-02-10 20:19:55'
$ws.Range('I39').Value = '10.6 mm'
$ws.Range('E40').Value = '2026-02-10 20:19:57'
$ws.Range('I40').Value = '13.4 mm'
$ws.Range('E41').Value = '2026-02-10 20:20:00'
$ws.Range('H41').Value = '''82%'
$ws.Range('J41').Value = '1004.6 hPa'
$ws.Range('O41').Value = '14.3 °C'
$ws.Range('E42').Value = '2026-02-10 20:20:02'
$ws.Range('E43').Value = '2026-02-10 20:20:04'
$ws.Range('O43').Value = '9.7 °C'
$ws.Range('E44').Value = '2026-02-10 20:20:07'
$ws.Range('H44').Value = '''96%'
$ws.Range('I44').Value = '28.3 mm'
$ws.Range('E45').Value = '2026-02-10 20:20:09'
$ws.Range('I45').Value = '37.2 mm'
$ws.Range('J45').Value = '1005.7 hPa'
$ws.Range('L45').Value = '27.4 km/h - 95º 19:45 TU'
$ws.Range('O45').Value = '6.5 °C'
$ws.Range('E46').Value = '2026-02-10 20:20:11'
$ws.Range('H46').Value = '''81%'
$ws.Range('O46').Value = '14.5 °C'
